# Auto-generated update of market-price-derived columns (H-N) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
  # Row 64
  $ws.Range("H64").Value = 3377.9
  $ws.Range("I64").Value = 3263.3333
  $ws.Range("J64").Value = 3427
  $ws.Range("K64").Value = 3263.3333
  $ws.Range("L64").Value = 3427
  $ws.Range("M64").Value = -3015.3333
  $ws.Range("N64").Value = -3923
  # Row 67
  $ws.Range("H67").Value = 3377.9
  $ws.Range("I67").Value = 3263.3333
  $ws.Range("J67").Value = 3427
  $ws.Range("K67").Value = 3263.3333
  $ws.Range("L67").Value = 3427
  $ws.Range("M67").Value = -2405.3333
  $ws.Range("N67").Value = -5143
  # Row 125
  $ws.Range("H125").Value = 474.92856
  $ws.Range("I125").Value = 391
  $ws.Range("K125").Value = 3519
  $ws.Range("M125").Value = -1059
  # Row 129
  $ws.Range("H129").Value = 141775.22
  $ws.Range("J129").Value = 164976.58
  $ws.Range("L129").Value = 494929.74
  $ws.Range("N129").Value = -504929.74
  # Row 137
  $ws.Range("H137").Value = 1721.4865
  $ws.Range("I137").Value = 1370.0741
  $ws.Range("J137").Value = 2670.3
  $ws.Range("K137").Value = 4110.2223
  $ws.Range("L137").Value = 8010.900000000001
  $ws.Range("M137").Value = -1560.2223
  $ws.Range("N137").Value = -13110.9

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
  # Row 74
  $ws.Range("H74").Value = 21740414
  $ws.Range("I74").Value = 26316328
  $ws.Range("J74").Value = 4824.875
  $ws.Range("K74").Value = 26316328
  $ws.Range("L74").Value = 4824.875
  $ws.Range("M74").Value = -26315454
  $ws.Range("N74").Value = -6572.875
  # Row 77
  $ws.Range("H77").Value = 21740414
  $ws.Range("I77").Value = 26316328
  $ws.Range("J77").Value = 4824.875
  $ws.Range("K77").Value = 131581640
  $ws.Range("L77").Value = 24124.375
  $ws.Range("M77").Value = -131577272
  $ws.Range("N77").Value = -32860.375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
  # Row 35
  $ws.Range("H35").Value = 24626.8
  $ws.Range("J35").Value = 24626.8
  $ws.Range("L35").Value = 24626.8
  $ws.Range("N35").Value = -25246.8
  # Row 94
  $ws.Range("H94").Value = 860.7273
  $ws.Range("I94").Value = 580.8889
  $ws.Range("J94").Value = 1196.5333
  $ws.Range("K94").Value = 580.8889
  $ws.Range("L94").Value = 1196.5333
  $ws.Range("M94").Value = -129.8889
  $ws.Range("N94").Value = -2098.5333
  # Row 135
  $ws.Range("H135").Value = 34243.75
  $ws.Range("J135").Value = 34243.75
  $ws.Range("L135").Value = 34243.75
  $ws.Range("N135").Value = -44383.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
  # Row 31
  $ws.Range("H31").Value = 4902.3
  $ws.Range("I31").Value = 3243.2
  $ws.Range("J31").Value = 5731.85
  $ws.Range("K31").Value = 3243.2
  $ws.Range("L31").Value = 5731.85
  $ws.Range("M31").Value = -2948.2
  $ws.Range("N31").Value = -6321.85
  # Row 34
  $ws.Range("H34").Value = 4902.3
  $ws.Range("I34").Value = 3243.2
  $ws.Range("J34").Value = 5731.85
  $ws.Range("K34").Value = 3243.2
  $ws.Range("L34").Value = 5731.85
  $ws.Range("M34").Value = -3041.2
  $ws.Range("N34").Value = -6135.85
  # Row 41
  $ws.Range("H41").Value = 19971.428
  $ws.Range("J41").Value = 19971.428
  $ws.Range("L41").Value = 19971.428
  $ws.Range("N41").Value = -20827.428
  # Row 122
  $ws.Range("H122").Value = 1313.8334
  $ws.Range("I122").Value = 1160.6923
  $ws.Range("J122").Value = 1712
  $ws.Range("K122").Value = 3482.0769
  $ws.Range("L122").Value = 5136
  $ws.Range("M122").Value = -1032.0769
  $ws.Range("N122").Value = -10036
  # Row 141
  $ws.Range("H141").Value = 22360.031
  $ws.Range("I141").Value = 8000
  $ws.Range("J141").Value = 22830.852
  $ws.Range("K141").Value = 8000
  $ws.Range("L141").Value = 22830.852
  $ws.Range("M141").Value = -2820
  $ws.Range("N141").Value = -33190.852

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
  # Row 4
  $ws.Range("H4").Value = 6000110
  $ws.Range("I4").Value = 225
  $ws.Range("J4").Value = 10000033
  $ws.Range("K4").Value = 675
  $ws.Range("L4").Value = 30000099
  $ws.Range("M4").Value = -563
  $ws.Range("N4").Value = -30000323
  # Row 117
  $ws.Range("H117").Value = 1213.3572
  $ws.Range("J117").Value = 1281
  $ws.Range("L117").Value = 3843
  $ws.Range("N117").Value = -10727
  # Row 131
  $ws.Range("H131").Value = 159504.02
  $ws.Range("J131").Value = 185996.72
  $ws.Range("L131").Value = 557990.16
  $ws.Range("N131").Value = -568070.16
  # Row 141
  $ws.Range("H141").Value = 4407
  $ws.Range("I141").Value = 1030
  $ws.Range("J141").Value = 5532.6665
  $ws.Range("K141").Value = 3090
  $ws.Range("L141").Value = 16597.9995
  $ws.Range("M141").Value = 2090
  $ws.Range("N141").Value = -26957.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
  # Row 46
  $ws.Range("H46").Value = 24826.166
  $ws.Range("I46").Value = 10020.5
  $ws.Range("J46").Value = 32229
  $ws.Range("K46").Value = 10020.5
  $ws.Range("L46").Value = 32229
  $ws.Range("M46").Value = -9864.5
  $ws.Range("N46").Value = -32541
  # Row 122
  $ws.Range("H122").Value = 3553.3333
  $ws.Range("I122").Value = 2551.0833
  $ws.Range("J122").Value = 4889.6665
  $ws.Range("K122").Value = 7653.249899999999
  $ws.Range("L122").Value = 14668.9995
  $ws.Range("M122").Value = -5203.249899999999
  $ws.Range("N122").Value = -19568.9995
  # Row 132
  $ws.Range("H132").Value = 15829.892
  $ws.Range("I132").Value = 1984.826
  $ws.Range("J132").Value = 38575.355
  $ws.Range("K132").Value = 5954.478
  $ws.Range("L132").Value = 115726.065
  $ws.Range("M132").Value = -3424.478
  $ws.Range("N132").Value = -120786.065

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
  # Row 7
  $ws.Range("H7").Value = 4470.8823
  $ws.Range("I7").Value = 4250
  $ws.Range("J7").Value = 4786.4287
  $ws.Range("K7").Value = 4250
  $ws.Range("L7").Value = 4786.4287
  $ws.Range("M7").Value = -4138
  $ws.Range("N7").Value = -5010.4287
  # Row 16
  $ws.Range("H16").Value = 570.8
  $ws.Range("I16").Value = 500.82352
  $ws.Range("K16").Value = 500.82352
  $ws.Range("M16").Value = -330.82352
  # Row 22
  $ws.Range("H22").Value = 2941.5
  $ws.Range("I22").Value = 4422.727
  $ws.Range("J22").Value = 1131.1111
  $ws.Range("K22").Value = 4422.727
  $ws.Range("L22").Value = 1131.1111
  $ws.Range("M22").Value = -4127.727
  $ws.Range("N22").Value = -1721.1111
  # Row 27
  $ws.Range("H27").Value = 2941.5
  $ws.Range("I27").Value = 4422.727
  $ws.Range("J27").Value = 1131.1111
  $ws.Range("K27").Value = 4422.727
  $ws.Range("L27").Value = 1131.1111
  $ws.Range("M27").Value = -4315.727
  $ws.Range("N27").Value = -1345.1111
  # Row 40
  $ws.Range("H40").Value = 3032.9
  $ws.Range("I40").Value = 2821.6086
  $ws.Range("K40").Value = 2821.6086
  $ws.Range("M40").Value = -2685.6086
  # Row 68
  $ws.Range("H68").Value = 2939.0715
  $ws.Range("I68").Value = 2950
  $ws.Range("J68").Value = 2924.5
  $ws.Range("K68").Value = 2950
  $ws.Range("L68").Value = 2924.5
  $ws.Range("M68").Value = -2201
  $ws.Range("N68").Value = -4422.5
  # Row 71
  $ws.Range("H71").Value = 2939.0715
  $ws.Range("I71").Value = 2950
  $ws.Range("J71").Value = 2924.5
  $ws.Range("K71").Value = 14750
  $ws.Range("L71").Value = 14622.5
  $ws.Range("M71").Value = -11006
  $ws.Range("N71").Value = -22110.5
  # Row 122
  $ws.Range("H122").Value = 579830.5
  $ws.Range("I122").Value = 894073.5
  $ws.Range("J122").Value = 3718.25
  $ws.Range("K122").Value = 2682220.5
  $ws.Range("L122").Value = 11154.75
  $ws.Range("M122").Value = -2679770.5
  $ws.Range("N122").Value = -16054.75
  # Row 126
  $ws.Range("H126").Value = 4470.8823
  $ws.Range("I126").Value = 4250
  $ws.Range("J126").Value = 4786.4287
  $ws.Range("K126").Value = 12750
  $ws.Range("L126").Value = 14359.2861
  $ws.Range("M126").Value = -10280
  $ws.Range("N126").Value = -19299.2861
  # Row 138
  $ws.Range("H138").Value = 45000
  $ws.Range("J138").Value = 45000
  $ws.Range("L138").Value = 45000
  $ws.Range("N138").Value = -55280

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
  # Row 113
  $ws.Range("H113").Value = 1347.5
  $ws.Range("I113").Value = 1544
  $ws.Range("J113").Value = 365
  $ws.Range("K113").Value = 4632
  $ws.Range("L113").Value = 1095
  $ws.Range("M113").Value = -2462
  $ws.Range("N113").Value = -5435
  # Row 126
  $ws.Range("H126").Value = 1857.8462
  $ws.Range("I126").Value = 1371.619
  $ws.Range("K126").Value = 4114.857
  $ws.Range("M126").Value = -1644.857
  # Row 136
  $ws.Range("H136").Value = 43015584
  $ws.Range("I136").Value = 68819480
  $ws.Range("K136").Value = 206458440
  $ws.Range("M136").Value = -206455890
